$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.144488096237183
$ws.Range("B1").Value = 2.522819042205811
$ws.Range("C1").Value = 4.051132678985596
$ws.Range("D1").Value = 3.545252084732056
$ws.Range("E1").Value = 1.225531101226807
